# Update odds values on the active worksheet (Jogos do Dia Betfair Back Lay)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("P2").Value = 2.12
$ws.Range("T2").Value = 1.85
$ws.Range("AC2").Value = 9

# Row 3 changes
$ws.Range("F3").Value = 1.38
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 11.5
$ws.Range("I3").Value = 16.5
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 5
$ws.Range("M3").Value = 1.09
$ws.Range("N3").Value = 2.84
$ws.Range("O3").Value = 1.43
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.28
$ws.Range("R3").Value = 1.22
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 2.62
$ws.Range("U3").Value = 1.49
$ws.Range("V3").Value = 1.07
$ws.Range("W3").Value = 3.05
$ws.Range("X3").Value = 13.5
$ws.Range("Y3").Value = 34
$ws.Range("AB3").Value = 5.8
$ws.Range("AC3").Value = 13.5
$ws.Range("AD3").Value = 65
$ws.Range("AF3").Value = 8
$ws.Range("AG3").Value = 14
$ws.Range("AH3").Value = 55
$ws.Range("AJ3").Value = 14
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 90
$ws.Range("AN3").Value = 13

# Row 4 changes
$ws.Range("H4").Value = 1.43
$ws.Range("V4").Value = 3.15
